$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.710855504857414
$ws.Range("D2").Value = 4.804792576929309
$ws.Range("E2").Value = 11.82195001753959
$ws.Range("F2").Value = 54.97338331793419
$ws.Range("G2").Value = 3.761966999061995
$ws.Range("I2").Value = 26.52218951118286
$ws.Range("J2").Value = 10.71073704869555
$ws.Range("K2").Value = 21.23899976025265
$ws.Range("B3").Value = 4.632236651735319
$ws.Range("D3").Value = 4.826469558780695
$ws.Range("E3").Value = 11.81653496596466
$ws.Range("F3").Value = 54.6311420076131
$ws.Range("G3").Value = 3.766536344042706
$ws.Range("I3").Value = 26.47044735959709
$ws.Range("J3").Value = 10.71851207039544
$ws.Range("K3").Value = 21.02325319522942
$ws.Range("B4").Value = 4.581987213235194
$ws.Range("D4").Value = 4.840719957687169
$ws.Range("E4").Value = 11.81606588957188
$ws.Range("F4").Value = 54.43228218873997
$ws.Range("G4").Value = 3.769483792932824
$ws.Range("I4").Value = 26.44133997575311
$ws.Range("J4").Value = 10.72518166284763
$ws.Range("K4").Value = 20.89680162844493
$ws.Range("B5").Value = 4.56102268477329
$ws.Range("D5").Value = 4.84676448661288
$ws.Range("E5").Value = 11.81659301440608
$ws.Range("F5").Value = 54.35412799357311
$ws.Range("G5").Value = 3.77072072426479
$ws.Range("I5").Value = 26.43014787860669
$ws.Range("J5").Value = 10.72837553597462
$ws.Range("K5").Value = 20.84684026074389
$ws.Range("B6").Value = 4.557512393035407
$ws.Range("D6").Value = 4.847782537760971
$ws.Range("E6").Value = 11.81672392141583
$ws.Range("F6").Value = 54.34132593597738
$ws.Range("G6").Value = 3.770928283784603
$ws.Range("I6").Value = 26.4283298571536
$ws.Range("J6").Value = 10.72893459917926
$ws.Range("K6").Value = 20.83864050821918
$ws.Range("B7").Value = 4.581706438825261
$ws.Range("D7").Value = 4.840800514063455
$ws.Range("E7").Value = 11.81607009037217
$ws.Range("F7").Value = 54.43121644322206
$ws.Range("G7").Value = 3.769500329377049
$ws.Range("I7").Value = 26.44118632531529
$ws.Range("J7").Value = 10.72522281043581
$ws.Range("K7").Value = 20.89612140978217
$ws.Range("B8").Value = 4.684163715696673
$ws.Range("D8").Value = 4.812072075109884
$ws.Range("E8").Value = 11.81949041689059
$ws.Range("F8").Value = 54.85306743637991
$ws.Range("G8").Value = 3.763513159449757
$ws.Range("I8").Value = 26.50379626594347
$ws.Range("J8").Value = 10.71302402076525
$ws.Range("K8").Value = 21.16339571854064
$ws.Range("B9").Value = 4.868999022917149
$ws.Range("D9").Value = 4.76316186905344
$ws.Range("E9").Value = 11.84883226585314
$ws.Range("F9").Value = 55.76752035965945
$ws.Range("G9").Value = 3.75289093338113
$ws.Range("I9").Value = 26.64769233640148
$ws.Range("J9").Value = 10.70417247980333
$ws.Range("K9").Value = 21.73266804167342
$ws.Range("B10").Value = 4.994505745172386
$ws.Range("D10").Value = 4.731703853410386
$ws.Range("E10").Value = 11.884129486133
$ws.Range("F10").Value = 56.48918826643218
$ws.Range("G10").Value = 3.745758983118855
$ws.Range("I10").Value = 26.7662590244058
$ws.Range("J10").Value = 10.70689014348933
$ws.Range("K10").Value = 22.17476473728582
$ws.Range("B11").Value = 5.049275398161893
$ws.Range("D11").Value = 4.71835404443023
$ws.Range("E11").Value = 11.90314852147962
$ws.Range("F11").Value = 56.82751701652799
$ws.Range("G11").Value = 3.742658349353739
$ws.Range("I11").Value = 26.82297259620997
$ws.Range("J11").Value = 10.71013318303855
$ws.Range("K11").Value = 22.38024445183805
$ws.Range("B12").Value = 5.069674024504153
$ws.Range("D12").Value = 4.713436105743908
$ws.Range("E12").Value = 11.91077380437134
$ws.Range("F12").Value = 56.95700445608222
$ws.Range("G12").Value = 3.741504725948892
$ws.Range("I12").Value = 26.84484560379695
$ws.Range("J12").Value = 10.71164987720873
$ws.Range("K12").Value = 22.45861200602134
$ws.Range("B13").Value = 5.065296111592104
$ws.Range("D13").Value = 4.71448917590834
$ws.Range("E13").Value = 11.90911279092116
$ws.Range("F13").Value = 56.92905730709009
$ws.Range("G13").Value = 3.741752269148522
$ws.Range("I13").Value = 26.84011725700287
$ws.Range("J13").Value = 10.71131039400549
$ws.Range("K13").Value = 22.44171052213714
$ws.Range("B14").Value = 5.050960488385071
$ws.Range("D14").Value = 4.717946694568822
$ws.Range("E14").Value = 11.90376739429683
$ws.Range("F14").Value = 56.83814300843688
$ws.Range("G14").Value = 3.742563029674261
$ws.Range("I14").Value = 26.82476416512769
$ws.Range("J14").Value = 10.71025217811362
$ws.Range("K14").Value = 22.38668104667669
$ws.Range("B15").Value = 5.042134831430753
$ws.Range("D15").Value = 4.720082387546522
$ws.Range("E15").Value = 11.90054820588938
$ws.Range("F15").Value = 56.78263150950193
$ws.Range("G15").Value = 3.743062311725061
$ws.Range("I15").Value = 26.81541154668237
$ws.Range("J15").Value = 10.70964157695967
$ws.Range("K15").Value = 22.35304430251898
$ws.Range("B16").Value = 4.990879021705368
$ws.Range("D16").Value = 4.732595556979611
$ws.Range("E16").Value = 11.88294592766651
$ws.Range("F16").Value = 56.46727344043308
$ws.Range("G16").Value = 3.745964496282229
$ws.Range("I16").Value = 26.76260833438068
$ws.Range("J16").Value = 10.70671859259126
$ws.Range("K16").Value = 22.16141796453332
$ws.Range("B17").Value = 4.958834674620152
$ws.Range("D17").Value = 4.740517450617893
$ws.Range("E17").Value = 11.87290435833737
$ws.Range("F17").Value = 56.27633096303574
$ws.Range("G17").Value = 3.747781598638979
$ws.Range("I17").Value = 26.73092488128718
$ws.Range("J17").Value = 10.70543950337504
$ws.Range("K17").Value = 22.04492955535656
$ws.Range("B18").Value = 4.940185269286541
$ws.Range("D18").Value = 4.745164377345765
$ws.Range("E18").Value = 11.86740772067365
$ws.Range("F18").Value = 56.16745538009253
$ws.Range("G18").Value = 3.748840284905276
$ws.Range("I18").Value = 26.71296296087275
$ws.Range("J18").Value = 10.70489269717622
$ws.Range("K18").Value = 21.97834391957953
$ws.Range("B19").Value = 4.933833630368929
$ws.Range("D19").Value = 4.746753308130067
$ws.Range("E19").Value = 11.86559465003064
$ws.Range("F19").Value = 56.13075725387881
$ws.Range("G19").Value = 3.749201067406575
$ws.Range("I19").Value = 26.70692635625824
$ws.Range("J19").Value = 10.70473999620659
$ws.Range("K19").Value = 21.95587267872667
$ws.Range("B20").Value = 4.962268502212623
$ws.Range("D20").Value = 4.739664794672082
$ws.Range("E20").Value = 11.87394444289288
$ws.Range("F20").Value = 56.29655937384507
$ws.Range("G20").Value = 3.747586765032136
$ws.Range("I20").Value = 26.73427058489944
$ws.Range("J20").Value = 10.70555611411417
$ws.Range("K20").Value = 22.05728746617497
$ws.Range("B21").Value = 5.055180530624025
$ws.Range("D21").Value = 4.716927416300257
$ws.Range("E21").Value = 11.90532600604531
$ws.Range("F21").Value = 56.86481017112492
$ws.Range("G21").Value = 3.742324333947564
$ws.Range("I21").Value = 26.82926299342958
$ws.Range("J21").Value = 10.71055516886418
$ws.Range("K21").Value = 22.40282999398517
$ws.Range("B22").Value = 5.113910665063433
$ws.Range("D22").Value = 4.702867410668921
$ws.Range("E22").Value = 11.92830071928747
$ws.Range("F22").Value = 57.24414236375193
$ws.Range("G22").Value = 3.739004564487054
$ws.Range("I22").Value = 26.89365744066908
$ws.Range("J22").Value = 10.71550458943867
$ws.Range("K22").Value = 22.63187460416844
$ws.Range("B23").Value = 5.082749933641848
$ws.Range("D23").Value = 4.710298542330063
$ws.Range("E23").Value = 11.91581416426529
$ws.Range("F23").Value = 57.04098353242762
$ws.Range("G23").Value = 3.740765499589061
$ws.Range("I23").Value = 26.85907836697559
$ws.Range("J23").Value = 10.71270908519816
$ws.Range("K23").Value = 22.50935881298339
$ws.Range("B24").Value = 4.960716774269452
$ws.Range("D24").Value = 4.740049992171752
$ws.Range("E24").Value = 11.87347335955435
$ws.Range("F24").Value = 56.2874113011063
$ws.Range("G24").Value = 3.747674805656657
$ws.Range("I24").Value = 26.73275720307894
$ws.Range("J24").Value = 10.70550280708498
$ws.Range("K24").Value = 22.05169924978816
$ws.Range("B25").Value = 4.820770464818086
$ws.Range("D25").Value = 4.775603635132723
$ws.Range("E25").Value = 11.8384755556592
$ws.Range("F25").Value = 55.51112927680872
$ws.Range("G25").Value = 3.755645779789803
$ws.Range("I25").Value = 26.60650481996252
$ws.Range("J25").Value = 10.70494964730955
$ws.Range("K25").Value = 21.57420176886598
